$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.882.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.985.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.00%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.553"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.53%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.39%  "

$ws.Range("E11").Value = "  +2.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.464.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.981.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.996"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.839.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("E25").Value = "  -4.80%  "

$ws.Range("E26").Value = "  -4.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  -0.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.109"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.65%  "

$ws.Range("E34").Value = "  +12.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.52%  "

$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("E40").Value = "  -4.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.31%  "

$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.09%  "

$ws.Range("E45").Value = "  -0.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.112.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.292.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.240"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.69%  "

$ws.Range("E51").Value = "  -0.92%  "
